$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.782.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.757.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "628.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.755.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.60%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.389.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.751.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.713.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "467.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.701"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("E25").Value = "  -6.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.903.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.178"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +19.90%  "
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.709.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.33%  "
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.954"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "155.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("E51").Value = "  -1.52%  "
